$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Class start date moved back one day (Sprint 1 now starts on class-begins day,
# no +1 offset in the Sprint-1 end-date formula below).
$ws.Range("B2").Value = 43702

# Sprint 1: duration 7 -> 14 days, formula no longer adds the extra day,
# and the "Labor Day" note moves from the Sprint 2 row up to the Sprint 1 row.
$ws.Range("B5").Formula = "=B2+D5"
$ws.Range("D5").Value = 14
$ws.Range("F5").Value = "Labor Day"

# Sprint 2 no longer carries the "Labor Day" note.
$ws.Range("F6").ClearContents()

# "Thanksgiving Break" note moves from the Sprint 8 row up to the Sprint 7 row,
# and is renamed with a trailing asterisk.
$ws.Range("F11").Value = "Thankgiving Break*"
$ws.Range("F12").ClearContents()

# A new blank formatted row appears below the old trailing blank row.
$ws.Rows.Item(16).RowHeight = 28

# Selection moved to A15.
$ws.Range("A15").Select()
